$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add older catch limit figures (B/D columns) for GB Cod, GOM Haddock,
# GOM Winter Flounder, SNE/MA Winter Flounder and SNE/MA Yellowtail Flounder.
$ws.Range("B13").Value = 10324
$ws.Range("D13").Value = 60729.411764705881

$ws.Range("B15").Value = 838
$ws.Range("D15").Value = 83800

$ws.Range("D17").Value = 62793.333333333336

# D19 switches from the unformatted style to the thousands-style already
# used by the other Target TAC cells (same look as D13/D15/D17/D25) -
# copy the format from a cell that already carries that style, then set
# the value.
$ws.Range("D17").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 61200

$ws.Range("D25").Value = 58816.666666666672

# Restore view: scroll to the top of the sheet and move the selection.
$ws.Range("J13").Select()
